$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Update summary figures at the top of the statement
# ---------------------------------------------------------------------------
# Valor Mora (total) changed
$ws.Range("E11").Value = 948000
# Cant. Periodos changed (2 -> 3 periods per worker now)
$ws.Range("F13").Value = 3

# ---------------------------------------------------------------------------
# 2) Make room for the new rows. The table currently has 8 data rows
#    (16-23, the last one carrying the "closing" thicker-border style).
#    The new layout needs 12 data rows (4 workers x 3 periods), so we need
#    4 additional rows, and the thicker closing border must end up on the
#    new last row (27) instead of the old one (23).
# ---------------------------------------------------------------------------
$ws.Rows("24:27").Insert()

# Preserve the special "closing" style (currently still sitting on row 23)
# by copying it down onto the new last row first ...
$ws.Range("B23:J23").Copy($ws.Range("B27:J27"))

# ... then convert row 23 back into a normal interior row, and give the
# freshly inserted rows 24-26 the same normal interior formatting.
$ws.Range("B16:J16").Copy($ws.Range("B23:J23"))
$ws.Range("B16:J16").Copy($ws.Range("B24:J24"))
$ws.Range("B16:J16").Copy($ws.Range("B25:J25"))
$ws.Range("B16:J16").Copy($ws.Range("B26:J26"))

# ---------------------------------------------------------------------------
# 3) Fill in the worker / period detail rows (16-27).
#    4 workers, each with 3 periods (2507, 2506, 2505), consecutively.
# ---------------------------------------------------------------------------
$workers = @(
    @{ Id = "45760135";   Name = "DAYLESTER CASARRUBIA LOPEZ";        Salario = 76000; Mora = 1900000 },
    @{ Id = "39280539";   Name = "ADRIANA LUCIA CRUZ LOPEZ";          Salario = 88000; Mora = 2200000 },
    @{ Id = "1010128177"; Name = "JEFERSON ANDRES PEREZ SALCEDO";     Salario = 76000; Mora = 1900000 },
    @{ Id = "3811397";    Name = "ANDRES ANTONIO CASARRUBIA LOPEZ";   Salario = 76000; Mora = 1900000 }
)
$periods = @("2507", "2506", "2505")

$row = 16
foreach ($worker in $workers) {
    foreach ($period in $periods) {
        $ws.Range("B$row").Value = "CC"
        $ws.Range("C$row").Value = $worker.Id
        $ws.Range("D$row").Value = $worker.Name
        $ws.Range("E$row").Value = $period
        $ws.Range("F$row").Value = $worker.Salario
        $ws.Range("G$row").Value = $worker.Mora
        $row = $row + 1
    }
}

Write-Output "done"
